$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 10

# Column A holds the date as plain text, like the existing rows above it.
# Prefix with a quote so Excel stores it as text instead of inferring a date
# serial number, then reset the cell style back to Normal so no stray
# number-format style lingers on the cell.
$ws.Cells.Item($row, 1).Value = "'12/04/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 14644.2
$ws.Cells.Item($row, 3).Value = 0.162034656703162
$ws.Cells.Item($row, 4).Value = 0.837965343296838
$ws.Cells.Item($row, 5).Value = -51.73
$ws.Cells.Item($row, 6).Value = -12.62
$ws.Cells.Item($row, 7).Value = -18195.42
$ws.Cells.Item($row, 8).Value = -59.72
$ws.Cells.Item($row, 9).Value = -429.98
$ws.Cells.Item($row, 10).Value = -15.34
